$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 02:52"

# --- Update Estados Unidos totals (row 4) ---
$ws.Cells.Item(4, 2).Value = 430271   # Casos totales
$ws.Cells.Item(4, 3).Value = 29936    # Nuevos casos
$ws.Cells.Item(4, 5).Value = 393177   # Casos activos
$ws.Cells.Item(4, 7).Value = 1897     # Casos criticos
$ws.Cells.Item(4, 8).Value = 14738    # Muertes

# --- Provincias/paises: Mayotte now has real data and is re-inserted into the
# --- sorted list right before "Islas Feroe" (previously it sat after "Kenia").
# --- Islas Feroe, Consejo Danes para los Refugiados and Kenia each keep their
# --- own stats but shift down one row to make room; the stale old Mayotte row
# --- is dropped.

# Row 114: was "Islas Feroe" -> becomes "Mayotte" with fresh data
$ws.Cells.Item(114, 1).Value = "Mayotte"
$ws.Cells.Item(114, 2).Value = 184
$ws.Cells.Item(114, 3).Value = 13
$ws.Cells.Item(114, 4).Value = 22
$ws.Cells.Item(114, 5).Value = 160
$ws.Cells.Item(114, 6).Value = 3
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 2

# Row 115: was "Consejo Danes para los Refugiados" -> becomes "Islas Feroe"
$ws.Cells.Item(115, 1).Value = "Islas Feroe"
$ws.Cells.Item(115, 2).Value = 184
$ws.Cells.Item(115, 3).Value = 0
$ws.Cells.Item(115, 4).Value = 131
$ws.Cells.Item(115, 5).Value = 53
$ws.Cells.Item(115, 6).Value = 1
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 0

# Row 116: was "Kenia" -> becomes "Consejo Danes para los Refugiados"
$ws.Cells.Item(116, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(116, 2).Value = 180
$ws.Cells.Item(116, 3).Value = 0
$ws.Cells.Item(116, 4).Value = 9
$ws.Cells.Item(116, 5).Value = 153
$ws.Cells.Item(116, 6).Value = 0
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 18

# Row 117: was "Mayotte" -> becomes "Kenia"
$ws.Cells.Item(117, 1).Value = "Kenia"
$ws.Cells.Item(117, 2).Value = 179
$ws.Cells.Item(117, 3).Value = 7
$ws.Cells.Item(117, 4).Value = 9
$ws.Cells.Item(117, 5).Value = 164
$ws.Cells.Item(117, 6).Value = 2
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(117, 8).Value = 6
